# Rename the inline picture shapes that live in the document's headers
# and footers:
#   - The Pearson logo picture (alt text / description
#     "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png")
#     appears twice (once per footer) and must be renamed from
#     "image2.png" to "image1.png".
#   - The BTEC logo picture (alt text / description "BTec_Logo-Orange")
#     appears once (in a header) and must be renamed from "image1.jpg"
#     to "image2.jpg".
#
# InlineShape objects obtained straight from a HeaderFooter.Range can be
# "stale" for the Name setter, so each shape reference is re-derived via
# its own Range right before the assignment.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $count = $hf.Range.InlineShapes.Count
            for ($i = 1; $i -le $count; $i++) {
                $shp = $hf.Range.InlineShapes($i)
                $shp = $shp.Range.InlineShapes(1)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                } elseif ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }

    foreach ($hf in $sec.Footers) {
        if ($hf.Exists) {
            $count = $hf.Range.InlineShapes.Count
            for ($i = 1; $i -le $count; $i++) {
                $shp = $hf.Range.InlineShapes($i)
                $shp = $shp.Range.InlineShapes(1)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                } elseif ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}

Write-Host "done"
